# Append a new trade record (row 13) to the BIIB noun-trade sheet,
# mirroring the structure of the existing rows (A:H).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13
$ws.Cells.Item($row, 1).Value = 9442.4
$ws.Cells.Item($row, 2).Value = 9482.23
$ws.Cells.Item($row, 3).Value = 311.98
$ws.Cells.Item($row, 4).Value = 310.67
$ws.Cells.Item($row, 5).Value = $false
$ws.Cells.Item($row, 6).Value = -0.42
$ws.Cells.Item($row, 7).Value = 42620.766111111108
$ws.Cells.Item($row, 8).Value = $false
